$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 468.2
$ws.Range("I15").Value = 468.2
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 1404.6
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -1235.6
$ws.Range("H17").Value = 2175384.8
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2175384.8
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6526154.399999999
$ws.Range("N17").Value = -6526490.399999999
$ws.Range("H32").Value = 833
$ws.Range("I32").Value = 462.25
$ws.Range("J32").Value = 920.2353000000001
$ws.Range("K32").Value = 462.25
$ws.Range("L32").Value = 920.2353000000001
$ws.Range("M32").Value = -136.25
$ws.Range("N32").Value = -1572.2353
$ws.Range("H64").Value = 3588.4314
$ws.Range("I64").Value = 3514.7222
$ws.Range("J64").Value = 3765.3333
$ws.Range("K64").Value = 3514.7222
$ws.Range("L64").Value = 3765.3333
$ws.Range("M64").Value = -3266.7222
$ws.Range("N64").Value = -4261.3333
$ws.Range("H67").Value = 3588.4314
$ws.Range("I67").Value = 3514.7222
$ws.Range("J67").Value = 3765.3333
$ws.Range("K67").Value = 3514.7222
$ws.Range("L67").Value = 3765.3333
$ws.Range("M67").Value = -2656.7222
$ws.Range("N67").Value = -5481.3333
$ws.Range("H74").Value = 4207.875
$ws.Range("I74").Value = 3334.3333
$ws.Range("J74").Value = 4732
$ws.Range("K74").Value = 3334.3333
$ws.Range("L74").Value = 4732
$ws.Range("M74").Value = -2398.3333
$ws.Range("N74").Value = -6604
$ws.Range("H77").Value = 4207.875
$ws.Range("I77").Value = 3334.3333
$ws.Range("J77").Value = 4732
$ws.Range("K77").Value = 16671.6665
$ws.Range("L77").Value = 23660
$ws.Range("M77").Value = -11991.6665
$ws.Range("N77").Value = -33020
$ws.Range("H88").Value = 2875.25
$ws.Range("I88").Value = 425
$ws.Range("J88").Value = 3365.3
$ws.Range("K88").Value = 425
$ws.Range("L88").Value = 3365.3
$ws.Range("M88").Value = -19
$ws.Range("N88").Value = -4177.3
$ws.Range("H91").Value = 2875.25
$ws.Range("I91").Value = 425
$ws.Range("J91").Value = 3365.3
$ws.Range("K91").Value = 425
$ws.Range("L91").Value = 3365.3
$ws.Range("M91").Value = 979
$ws.Range("N91").Value = -6173.3
$ws.Range("H129").Value = 198896.3
$ws.Range("I129").Value = 612.8333
$ws.Range("J129").Value = 225334.1
$ws.Range("K129").Value = 1838.4999
$ws.Range("L129").Value = 676002.3
$ws.Range("M129").Value = 3161.5001
$ws.Range("N129").Value = -686002.3
$ws.Range("H138").Value = 3022.68
$ws.Range("I138").Value = 1588.9791
$ws.Range("J138").Value = 4346.096
$ws.Range("K138").Value = 4766.9373
$ws.Range("L138").Value = 13038.288
$ws.Range("M138").Value = 373.0627000000004
$ws.Range("N138").Value = -23318.288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -384
$ws.Range("N4").Value = $null
$ws.Range("H6").Value = 41876.5
$ws.Range("I6").Value = 41876.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 41876.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -41703.5
$ws.Range("N6").Value = $null
$ws.Range("H9").Value = 14000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 14000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 14000
$ws.Range("N9").Value = -14340
$ws.Range("H20").Value = 14000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 14000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 14000
$ws.Range("N20").Value = -14540
$ws.Range("H23").Value = 73879.75
$ws.Range("I23").Value = 80006
$ws.Range("J23").Value = 67753.5
$ws.Range("K23").Value = 80006
$ws.Range("L23").Value = 67753.5
$ws.Range("M23").Value = -79747
$ws.Range("N23").Value = -68271.5
$ws.Range("H37").Value = 12222.223
$ws.Range("I37").Value = 12222.223
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 12222.223
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -11949.223
$ws.Range("H61").Value = 2443.074
$ws.Range("I61").Value = 2556.8462
$ws.Range("J61").Value = 2337.4285
$ws.Range("K61").Value = 2556.8462
$ws.Range("L61").Value = 2337.4285
$ws.Range("M61").Value = -2344.8462
$ws.Range("N61").Value = -2761.4285
$ws.Range("H136").Value = 2443.074
$ws.Range("I136").Value = 2556.8462
$ws.Range("J136").Value = 2337.4285
$ws.Range("K136").Value = 7670.5386
$ws.Range("L136").Value = 7012.2855
$ws.Range("M136").Value = -5120.5386
$ws.Range("N136").Value = -12112.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 414.7143
$ws.Range("I22").Value = 320.6
$ws.Range("J22").Value = 650
$ws.Range("K22").Value = 320.6
$ws.Range("L22").Value = 650
$ws.Range("M22").Value = -147.6
$ws.Range("N22").Value = -996
$ws.Range("H68").Value = 50000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 50000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 50000
$ws.Range("N68").Value = -51622
$ws.Range("H71").Value = 50000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 50000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 150000
$ws.Range("N71").Value = -158112
$ws.Range("H82").Value = 8551.4
$ws.Range("I82").Value = 1989.25
$ws.Range("J82").Value = 34800
$ws.Range("K82").Value = 1989.25
$ws.Range("L82").Value = 34800
$ws.Range("M82").Value = -1606.25
$ws.Range("H85").Value = 8551.4
$ws.Range("I85").Value = 1989.25
$ws.Range("J85").Value = 34800
$ws.Range("K85").Value = 1989.25
$ws.Range("L85").Value = 34800
$ws.Range("M85").Value = -663.25
$ws.Range("H86").Value = 1889.2979
$ws.Range("I86").Value = 1592.72
$ws.Range("J86").Value = 2226.318
$ws.Range("K86").Value = 1592.72
$ws.Range("L86").Value = 2226.318
$ws.Range("M86").Value = -469.72
$ws.Range("H89").Value = 1889.2979
$ws.Range("I89").Value = 1592.72
$ws.Range("J89").Value = 2226.318
$ws.Range("K89").Value = 7963.6
$ws.Range("L89").Value = 11131.59
$ws.Range("M89").Value = -2347.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1174.1
$ws.Range("I7").Value = 1293.4445
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 1293.4445
$ws.Range("L7").Value = 100
$ws.Range("M7").Value = -1180.4445
$ws.Range("H22").Value = 1543.3334
$ws.Range("I22").Value = 2173.3333
$ws.Range("J22").Value = 283.33334
$ws.Range("K22").Value = 2173.3333
$ws.Range("L22").Value = 283.33334
$ws.Range("M22").Value = -1823.3333
$ws.Range("N22").Value = -983.33334
$ws.Range("H31").Value = 3402.9277
$ws.Range("I31").Value = 1405.4407
$ws.Range("J31").Value = 8313.416999999999
$ws.Range("K31").Value = 1405.4407
$ws.Range("L31").Value = 8313.416999999999
$ws.Range("M31").Value = -1110.4407
$ws.Range("N31").Value = -8903.416999999999
$ws.Range("H34").Value = 3402.9277
$ws.Range("I34").Value = 1405.4407
$ws.Range("J34").Value = 8313.416999999999
$ws.Range("K34").Value = 1405.4407
$ws.Range("L34").Value = 8313.416999999999
$ws.Range("M34").Value = -1203.4407
$ws.Range("N34").Value = -8717.416999999999
$ws.Range("H62").Value = 4374.067
$ws.Range("I62").Value = 3245
$ws.Range("J62").Value = 6067.6665
$ws.Range("K62").Value = 3245
$ws.Range("L62").Value = 6067.6665
$ws.Range("M62").Value = -2621
$ws.Range("N62").Value = -7315.6665
$ws.Range("H65").Value = 4374.067
$ws.Range("I65").Value = 3245
$ws.Range("J65").Value = 6067.6665
$ws.Range("K65").Value = 16225
$ws.Range("L65").Value = 30338.3325
$ws.Range("M65").Value = -13105
$ws.Range("N65").Value = -36578.3325

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 845.09
$ws.Range("I131").Value = 556.3333
$ws.Range("J131").Value = 854.0205999999999
$ws.Range("K131").Value = 1668.9999
$ws.Range("L131").Value = 2562.0618
$ws.Range("M131").Value = 3371.0001
$ws.Range("N131").Value = -12642.0618

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2235.889
$ws.Range("I68").Value = 1608.6
$ws.Range("J68").Value = 3020
$ws.Range("K68").Value = 1608.6
$ws.Range("L68").Value = 3020
$ws.Range("M68").Value = -859.5999999999999
$ws.Range("N68").Value = -4518
$ws.Range("H71").Value = 2235.889
$ws.Range("I71").Value = 1608.6
$ws.Range("J71").Value = 3020
$ws.Range("K71").Value = 8043
$ws.Range("L71").Value = 15100
$ws.Range("M71").Value = -4299
$ws.Range("N71").Value = -22588
$ws.Range("H132").Value = 5634.927
$ws.Range("I132").Value = 5944.5
$ws.Range("J132").Value = 4886.7915
$ws.Range("K132").Value = 17833.5
$ws.Range("L132").Value = 14660.3745
$ws.Range("M132").Value = -15303.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5210426.5
$ws.Range("I136").Value = 8547473
$ws.Range("J136").Value = 4634.8
$ws.Range("K136").Value = 25642419
$ws.Range("L136").Value = 13904.4
$ws.Range("M136").Value = -25639869
$ws.Range("N136").Value = -19004.4
